$wb = $excel.ActiveWorkbook

# --- Sheet 1: "All Published Values" -- append new published-rate row ---
$ws1 = $wb.Worksheets.Item("All Published Values")

$newRow = 10
$ws1.Cells.Item($newRow, 1).Value  = "'2026-01-02"
$ws1.Cells.Item($newRow, 2).Value  = "2026-01-02 19:48:09"
$ws1.Cells.Item($newRow, 3).Value  = "'697.85"
$ws1.Cells.Item($newRow, 4).Value  = "'697.85"
$ws1.Cells.Item($newRow, 5).Value  = "'700.79"
$ws1.Cells.Item($newRow, 6).Value  = "'700.79"
$ws1.Cells.Item($newRow, 7).Value  = "'702.88"
$ws1.Cells.Item($newRow, 8).Value  = "2026/01/02 19:48:09"
$ws1.Cells.Item($newRow, 9).Value  = "2026-01-02 11:51:31"
$ws1.Cells.Item($newRow, 10).Value = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"

# Re-point the AutoFilter so its range covers the freshly added row.
# (Toggle off first -- Range.AutoFilter() flips the existing filter off
# when one is already active on the sheet, so this sequence lands "on"
# with the new $A$1:$J$10 range instead of silently removing it.)
$ws1.AutoFilterMode = $false
$null = $ws1.Range("A1:J10").AutoFilter()

# The hidden _xlnm._FilterDatabase defined name tracks the autofilter
# range too, but isn't resynced automatically -- update it explicitly.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "All Published Values!_FilterDatabase") {
        $n.RefersTo = "='All Published Values'!`$A`$1:`$J`$10"
    }
}

# --- Sheet 2: "Daily Summary" -- bump today's publish count ---
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Cells.Item(4, 2).Value = 9
